$wb = $excel.ActiveWorkbook

# --- loads sheet: insert new columns (s_base_mva, v_nom_pu shift, g_shunt_pu, b_shunt_pu) ---
$wsLoads = $wb.Worksheets.Item("loads")

# Header row
$wsLoads.Range("A1").Value = "name"
$wsLoads.Range("B1").Value = "v_nom_kv"
$wsLoads.Range("C1").Value = "s_base_mva"
$wsLoads.Range("D1").Value = "v_nom_pu"
$wsLoads.Range("E1").Value = "p_nom_mw"
$wsLoads.Range("F1").Value = "q_nom_mvar"
$wsLoads.Range("G1").Value = "bus_idx"
$wsLoads.Range("H1").Value = "g_shunt_pu"
$wsLoads.Range("I1").Value = "b_shunt_pu"

# Data row
$wsLoads.Range("A2").Value = "Load 1"
$wsLoads.Range("B2").Value = 22
$wsLoads.Range("C2").Value = 100
$wsLoads.Range("D2").Value = 1
$wsLoads.Range("E2").Value = 10
$wsLoads.Range("F2").Value = 10
$wsLoads.Range("G2").Value = 2
$wsLoads.Range("H2").Value = 0
$wsLoads.Range("I2").Value = 0

[void]$wsLoads.Range("J2").Select()

# --- trafos sheet: add new columns (idx_hv, idx_lv, tap_pos, tap_change, tap_min, tap_max, v_base_kV) ---
$wsTrafos = $wb.Worksheets.Item("trafos")

$wsTrafos.Range("A1").Value = "name"
$wsTrafos.Range("B1").Value = "S_nom"
$wsTrafos.Range("C1").Value = "V_hv_kV"
$wsTrafos.Range("D1").Value = "V_lv_kV"
$wsTrafos.Range("F1").Value = "V_SCH_pu"
$wsTrafos.Range("G1").Value = "P_Cu_pu"
$wsTrafos.Range("H1").Value = "I_E_pu"
$wsTrafos.Range("I1").Value = "P_Fe_pu"
$wsTrafos.Range("J1").Value = "idx_hv"
$wsTrafos.Range("K1").Value = "idx_lv"
$wsTrafos.Range("L1").Value = "tap_pos"
$wsTrafos.Range("M1").Value = "tap_change"
$wsTrafos.Range("N1").Value = "tap_min"
$wsTrafos.Range("O1").Value = "tap_max"
$wsTrafos.Range("E1").Value = "v_base_kV"

[void]$wsTrafos.Range("E1").Select()
[void]$wsTrafos.Activate()
